$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.641.02"
$ws.Range("E2").Value = "  +4.57%  "
$ws.Range("D3").Value = "2.726.00"
$ws.Range("E3").Value = "  +2.78%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.47%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.995"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("E8").Value = "  +1.35%  "
$ws.Range("D9").Value = "2.751.36"
$ws.Range("E9").Value = "  +2.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.70"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.113"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.98%  "
$ws.Range("B12").Value = "Cardano"
$ws.Range("C12").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.390"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.64%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.162"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.06%  "
$ws.Range("D14").Value = "3.210.47"
$ws.Range("E14").Value = "  +2.93%  "
$ws.Range("E15").Value = "  +1.94%  "
$ws.Range("D16").Value = "63.542.88"
$ws.Range("E16").Value = "  +4.45%  "
$ws.Range("E17").Value = "  +5.98%  "
$ws.Range("D18").Value = "2.743.58"
$ws.Range("E18").Value = "  +3.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.85%  "
$ws.Range("E20").Value = "  +2.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "360.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.95"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.28%  "
$ws.Range("E24").Value = "  +0.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.75%  "
$ws.Range("E26").Value = "  +4.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.53"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.997"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("D29").Value = "0.0₃0909"
$ws.Range("E29").Value = "  +11.88%  "
$ws.Range("E30").Value = "  -0.94%  "
$ws.Range("E31").Value = "  +6.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "171.14"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.50%  "
$ws.Range("E33").Value = "  +13.02%  "
$ws.Range("E34").Value = "  -0.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "20.50"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.85%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.78"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.56%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.43"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.79%  "
$ws.Range("E38").Value = "  +9.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.995"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +13.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "344.52"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.31"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.59"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.76"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.88"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0589"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "139.35"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.41%  "
$ws.Range("E48").Value = "  +4.45%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0255"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.31%  "
$ws.Range("E50").Value = "  +0.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.995"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.25%  "
